$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-05 Tuesday", "2025-08-06 Wednesday"),
    @("20+11=31", "94-44=50"),
    @("90-61=29", "63-13=50"),
    @("78-71=7", "73+25=98"),
    @("92-11=81", "47+15=62"),
    @("62-7=55", "66-8=58"),
    @("36-9=27", "32+56=88"),
    @("99-64=35", "91-51=40"),
    @("88-10=78", "65+23=88"),
    @("90+4=94", "95-43=52"),
    @("23+11=34", "93-75=18"),
    @("54+34=88", "66-18=48"),
    @("3+91=94", "79-78=1"),
    @("99-86=13", "94-15=79"),
    @("31+7=38", "9+35=44"),
    @("48-29=19", "15+30=45"),
    @("68-13=55", "82+4=86"),
    @("3+73=76", "54+19=73"),
    @("82+7=89", "43-11=32"),
    @("78+12=90", "37-35=2"),
    @("89+10=99", "17+52=69"),
    @("62-28=34", "5+77=82"),
    @("78-10=68", "37+27=64"),
    @("31-30=1", "27+36=63"),
    @("67-49=18", "43+44=87"),
    @("54-30=24", "49-2=47"),
    @("73-27=46", "18-16=2"),
    @("67-59=8", "70-44=26"),
    @("58+1=59", "98-65=33"),
    @("19+78=97", "94-87=7"),
    @("91-47=44", "89-51=38"),
    @("1+33=34", "77-6=71"),
    @("82-50=32", "57-36=21"),
    @("69-59=10", "31+38=69"),
    @("60+9=69", "47+0=47"),
    @("28+23=51", "18-6=12"),
    @("82-6=76", "74-64=10"),
    @("20+27=47", "96-31=65"),
    @("55-19=36", "24+56=80"),
    @("87-40=47", "23+72=95"),
    @("45+54=99", "68-27=41"),
    @("25+30=55", "9+44=53"),
    @("0+82=82", "22-21=1"),
    @("86-0=86", "14+0=14"),
    @("8+18=26", "0+32=32"),
    @("51-38=13", "60+12=72"),
    @("72-21=51", "25-19=6"),
    @("84-77=7", "49+4=53"),
    @("93-37=56", "70-7=63"),
    @("36-27=9", "17+60=77"),
    @("16-14=2", "47-12=35"),
    @("30+20=50", "79-6=73"),
    @("69-10=59", "58+16=74"),
    @("22-18=4", "24+69=93"),
    @("8+35=43", "76+11=87"),
    @("41+14=55", "62-59=3"),
    @("11+45=56", "85+7=92"),
    @("63-22=41", "48-46=2"),
    @("5+92=97", "0+81=81"),
    @("77-25=52", "46-5=41"),
    @("93-47=46", "18+18=36"),
    @("44+18=62", "0+52=52"),
    @("45+52=97", "41-20=21"),
    @("19-15=4", "79-49=30"),
    @("93-4=89", "11+58=69"),
    @("57-19=38", "17+21=38"),
    @("62-47=15", "39+35=74"),
    @("57+11=68", "17+72=89"),
    @("72-20=52", "3+49=52"),
    @("81-62=19", "74+11=85"),
    @("25+34=59", "9+39=48"),
    @("10+71=81", "71-0=71"),
    @("57-39=18", "56+0=56"),
    @("37+36=73", "43+18=61"),
    @("64-57=7", "82+16=98"),
    @("17+2=19", "3+26=29"),
    @("34-1=33", "79-73=6"),
    @("25+62=87", "25+22=47"),
    @("19-12=7", "37+25=62"),
    @("38+9=47", "87-6=81"),
    @("91-75=16", "25-0=25"),
    @("56-1=55", "14+22=36"),
    @("91+2=93", "37-15=22"),
    @("14+71=85", "10+8=18"),
    @("70-69=1", "63+13=76"),
    @("33+49=82", "98-26=72"),
    @("11+41=52", "38+37=75"),
    @("55+12=67", "43+49=92"),
    @("7+38=45", "37+59=96"),
    @("19+12=31", "1+24=25"),
    @("94-69=25", "16+15=31"),
    @("14+32=46", "42-36=6"),
    @("58-26=32", "61-32=29"),
    @("85-8=77", "4+48=52"),
    @("45-10=35", "79-5=74"),
    @("16+11=27", "9-1=8"),
    @("33+2=35", "35+31=66"),
    @("85+8=93", "37+38=75"),
    @("91-11=80", "96-2=94"),
    @("83-43=40", "97-42=55"),
    @("97-64=33", "77-45=32"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}
